$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Merge the two runs that make up "SAT Dec 9 ... 2017" into a single
#    run (the split between "SAT Dec 9" and " 11:11:11 PST 2017" goes
#    away).
# ---------------------------------------------------------------------
$null = $d.Content.Find.Execute("SAT Dec 9 11:11:11 PST 2017", $false, $false, $false, $false, $false, $true, 1, $false, "SAT Dec 9 11:11:11 PST 2017", 2)

# ---------------------------------------------------------------------
# 2) Append a brand new "purchase" record after the last one in the
#    document (SAT Dec 9 .../Amount Received mode - CASH block), i.e.
#    a "18/12/2017 MAMATHA CHICK IN" style new entry dated
#    "SAT Dec 16 12:41:27 PST 2017".
# ---------------------------------------------------------------------

# Locate the very last paragraph that still carries text (the existing
# last record ends in "Amount Received mode" / "- CASH", immediately
# followed by a run of blank PlainText paragraphs).
$lastTextIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    $clean = $t -replace "`r", ""
    if ($clean -ne "") {
        $lastTextIndex = $i
    }
}

$anchor = $d.Paragraphs.Item($lastTextIndex)

function New-Row {
    param($afterParagraph, $label, $tabCount, $value, $fmt)

    $null = $afterParagraph.Range.InsertParagraphAfter()
    $idx = $afterParagraph.Index + 1
    $target = $d.Paragraphs.Item($idx)

    $tabs = "`t" * $tabCount
    $target.Range.Text = "$label$tabs$value"

    if ($fmt -eq "red") {
        $target.Range.Font.Color = 255
    } elseif ($fmt -eq "bold") {
        $target.Range.Font.Bold = 1
    }

    return $target
}

function New-BlankRow {
    param($afterParagraph)
    $null = $afterParagraph.Range.InsertParagraphAfter()
    $idx = $afterParagraph.Index + 1
    return $d.Paragraphs.Item($idx)
}

function New-DateRow {
    param($afterParagraph, $datePart, $timePart)

    $null = $afterParagraph.Range.InsertParagraphAfter()
    $idx = $afterParagraph.Index + 1
    $target = $d.Paragraphs.Item($idx)

    $target.Range.Text = "$datePart$timePart"

    # Force the leading "date" text and the trailing " time" text into
    # two separate runs (matching the source document's convention for
    # date/time stamp paragraphs) by toggling a character property on
    # just the first part and then reverting it.
    $s = $target.Range.Start
    $sub = $d.Range($s, $s + $datePart.Length)
    $sub.Font.Bold = 1
    $sub.Font.Bold = 0

    return $target
}

# Blank separator line straight after the old last record.
$anchor = New-BlankRow $anchor

# New record header line.
$anchor = New-DateRow $anchor "SAT Dec 16" " 12:41:27 PST 2017"

$anchor = New-Row $anchor "Person Name" 4 "- MAU" ""
$anchor = New-Row $anchor "Bill number" 4 "- 86" ""
$anchor = New-Row $anchor "---------------------------------------------------------------" 0 "" ""
$anchor = New-Row $anchor "Item Name" 4 "- CARROT" ""
$anchor = New-Row $anchor "Number of Pockets" 3 "- 1" ""
$anchor = New-Row $anchor "Number of KGs" 3 "- 86" ""
$anchor = New-Row $anchor "Rate" 5 "- 40" ""
$anchor = New-Row $anchor "Total Price" 4 "- 3440.0" ""
$anchor = New-Row $anchor "Amount Received" 3 "- 4000" "red"
$anchor = New-Row $anchor "Amount balance" 3 "- 38431.0" "bold"
$anchor = New-Row $anchor "Amount Received mode" 2 "- CASH" ""

# Two trailing blank PlainText paragraphs after the new record.
$anchor = New-BlankRow $anchor
$anchor = New-BlankRow $anchor

Write-Host "Paragraph count now:" $d.Paragraphs.Count
